# Auto-generated edit script: updates crypto price/volume table per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "79.907.44"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +4.48%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.226.00"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +5.62%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  -0.07%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "206.98"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +2.43%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "642.86"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +2.84%  "

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -0.05%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.241"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +16.02%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.585"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +5.99%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.216.58"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +5.39%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.581"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +31.84%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.166"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  +3.01%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.55"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +7.06%  "

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  +19.80%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.805.23"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +5.22%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "32.10"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +8.90%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "79.771.84"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +4.41%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.209.28"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +4.93%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "14.58"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +7.15%  "

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +31.07%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "9.28"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +2.19%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "433.01"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +15.28%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.15"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +17.69%  "

$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.377.66"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +5.25%  "

$ws.Range("B25").Value = "Aptos"
$ws.Range("C25").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "11.32"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +14.12%  "

$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "4.77"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +7.95%  "

$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "77.09"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  +4.69%  "

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +0.14%  "

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +7.46%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.10"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  +9.33%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  +0.16%  "

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +5.60%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "529.96"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +4.52%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.01"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +2.71%  "

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +22.41%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "23.27"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +11.43%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.121"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +13.74%  "

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -0.03%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "164.88"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +1.33%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "20.04"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  +0.02%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "193.69"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +0.42%  "

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  +0.00%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.57"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +7.27%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.825"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +4.60%  "

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +8.31%  "

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  +4.04%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "43.44"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +2.96%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.649"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +5.65%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "26.02"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +15.67%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.56"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +3.56%  "
